$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new date columns before column B (B:D), pushing the old
# B:E date columns right to E:H.
$ws.Columns("B:D").Insert()

# Keep the date columns the same fixed width as before.
$ws.Columns("C:H").ColumnWidth = 7.14

# Header row: two new dates, newest on the left.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Default every new cell in the new date columns to "UN" (unchanged).
$ws.Range("B2:D27").Value = "UN"

# Morgan Stanley (row 7): price-target cut on 6/25/2018.
$ws.Range("C7:D7").Value = "6/25/2018,Lowers Target,Equal Weight -> Equal Weight,$26.00 -> $23.00"
$ws.Range("C7:D7").Interior.ColorIndex = 45

# BidaskClub (row 22): downgrade on 6/26/2018.
$ws.Range("B22:D22").Value = "6/26/2018,Downgrades,Buy -> Hold,"
$ws.Range("D22").Interior.ColorIndex = 45

# New coverage added for this ticker.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
